$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$oldText = $ws.Range("C2").Value2
$newText = $oldText -replace [regex]::Escape("-->(:arm)-->"), "-->(a:arm)-->"

$ws.Range("C2").Value2 = $newText
$ws.Range("C3").Value2 = $newText

$ws.Range("B3").Select()
